# [2023-11-22] Spring boot 기본 end
#
# On slide 13 (TextBox 5) two paragraphs had their text accidentally split
# across multiple runs with identical formatting. Re-merge them back into
# single runs, exactly as a user retyping/joining the text in PowerPoint
# would produce:
#   - Paragraph "1. 복합적인 ..."       : runs "1" + ". "            -> "1. "
#   - Paragraph " 1) @interface 생성"   : runs "1) " + "@" + "interface " -> "1) @interface "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# --- Edit 1: merge "1" + ". " -> "1. " at the start of paragraph 1 ---
$para1 = $tr.Paragraphs(1, 1)
$tr.Characters($para1.Start, 3).Text = "1. "

# --- Edit 2: merge "1) " + "@" + "interface " -> "1) @interface " in paragraph 7 ---
# Paragraph 7 text is " 1) @interface 생성"; skip the leading space run (1 char)
# then the next 14 characters are "1) @interface ".
$para7 = $tr.Paragraphs(7, 1)
$tr.Characters($para7.Start + 1, 14).Text = "1) @interface "
